$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "readme" sheet: reorder the Author / JobNo / Date columns (B, C, D)
#    to JobNo / Date / Author, keeping header text and row data together.
# ---------------------------------------------------------------------
$readme = $wb.Worksheets.Item("readme")

# Header row (row 1)
$readme.Range("B1").Value = "JobNo"
$readme.Range("C1").Value = "Date"
$readme.Range("D1").Value = "Author"

# Data rows (rows 2-12): rotate the existing cell contents (rather than
# re-typing the values) so that the numeric-looking text "20220303"
# keeps its original text type/format instead of being re-interpreted
# as a number.
$readme.Range("B2:B12").Copy()
$readme.Range("F2:F12").PasteSpecial(-4163)

$readme.Range("C2:C12").Copy()
$readme.Range("B2:B12").PasteSpecial(-4163)

$readme.Range("D2:D12").Copy()
$readme.Range("C2:C12").PasteSpecial(-4163)

$readme.Range("F2:F12").Copy()
$readme.Range("D2:D12").PasteSpecial(-4163)

$readme.Range("F2:F12").Clear()

# ---------------------------------------------------------------------
# 2) "Project Information" sheet: update the analysis timestamp.
# ---------------------------------------------------------------------
$projInfo = $wb.Worksheets.Item("Project Information")
$projInfo.Range("B11").Value = "2022-03-03 15:33:43.320121"

# ---------------------------------------------------------------------
# 3) "Criterion Definitions" sheet: rename Criterion 2 definition text.
# ---------------------------------------------------------------------
$critDefs = $wb.Worksheets.Item("Criterion Definitions")
$critDefs.Range("A3").Value = "Criterion 2 (Max Daily Weight)"

# ---------------------------------------------------------------------
# 4) Results sheets: rename the "Criterion 2" column header (column F)
#    on each of the air-speed results sheets.
# ---------------------------------------------------------------------
$resultSheets = @(
    "Results, Air Speed 0.1",
    "Results, Air Speed 0.15",
    "Results, Air Speed 0.2",
    "Results, Air Speed 0.3",
    "Results, Air Speed 0.4",
    "Results, Air Speed 0.5",
    "Results, Air Speed 0.6",
    "Results, Air Speed 0.7",
    "Results, Air Speed 0.8"
)

foreach ($sheetName in $resultSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F1").Value = "Criterion 2 (Max Daily Weight)"
}
